# Bug report tracker: add bug #11 describing the Location class refactor
# (currentLocation no longer prints hostility after the constructor change).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "currentLocation method not printing hostility after location constructor refactor"
$ws.Range("C12").Value = "ricky"
# Write E12 before D12 so the shared-string table gets the same new-entry
# ordering (index 40 = Fix text, index 41 = Cause text) as the source file.
$ws.Range("E12").Value = "Refactor Location class to include hostility parameter "
$ws.Range("D12").Value = "locationclass wasn’t set up properly"
$ws.Range("F12").Value = "fixed"

# Column B needed to widen to fit the new, longer description text.
$ws.Columns.Item(2).ColumnWidth = 72

# Final cursor position left on D9 after editing the new row.
$ws.Range("D9").Select()
